# Applies the scheduled-runner Sheets update to the Marilith Profits workbook.
# For each affected leve row, refresh the market-price / profit columns
# (currentAveragePrice.. / LeveProfit..) to the latest recalculated figures.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 86
$ws.Range("H86").Value = 10597
$ws.Range("I86").Value = 4580.8335
$ws.Range("K86").Value = 4580.8335
$ws.Range("M86").Value = -3457.8335

# Row 89
$ws.Range("H89").Value = 10597
$ws.Range("I89").Value = 4580.8335
$ws.Range("K89").Value = 22904.1675
$ws.Range("M89").Value = -17288.1675

$ws = $wb.Worksheets.Item("BSM")
# Row 26
$ws.Range("H26").Value = 13333
$ws.Range("I26").Value = 13333
$ws.Range("K26").Value = 13333
$ws.Range("M26").Value = -13041

$ws = $wb.Worksheets.Item("CRP")
# Row 13
$ws.Range("H13").Value = 0
$ws.Range("I13").Value = 0
$ws.Range("J13").Value = 0
$ws.Range("K13").Value = 0
$ws.Range("L13").Value = 0
$ws.Range("M13").ClearContents()
$ws.Range("N13").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
# Row 11
$ws.Range("H11").Value = 652.2143
$ws.Range("I11").Value = 535.2857
$ws.Range("J11").Value = 769.1429000000001
$ws.Range("K11").Value = 1605.8571
$ws.Range("L11").Value = 2307.4287
$ws.Range("M11").Value = -1465.8571
$ws.Range("N11").Value = -2587.4287

# Row 23
$ws.Range("H23").Value = 711.8570999999999
$ws.Range("J23").Value = 621.6
$ws.Range("L23").Value = 1864.8
$ws.Range("N23").Value = -2334.8

# Row 63
$ws.Range("H63").Value = 2000
$ws.Range("I63").Value = 2000
$ws.Range("K63").Value = 6000
$ws.Range("M63").Value = -5251

# Row 64
$ws.Range("H64").Value = 1500
$ws.Range("I64").Value = 500
$ws.Range("K64").Value = 1500
$ws.Range("M64").Value = -1230

# Row 66
$ws.Range("H66").Value = 2000
$ws.Range("I66").Value = 2000
$ws.Range("K66").Value = 18000
$ws.Range("M66").Value = -14256

# Row 67
$ws.Range("H67").Value = 1500
$ws.Range("I67").Value = 500
$ws.Range("K67").Value = 1500
$ws.Range("M67").Value = -564

# Row 108
$ws.Range("H108").Value = 792.6667
$ws.Range("I108").Value = 792.6667
$ws.Range("K108").Value = 2378.0001
$ws.Range("M108").Value = 501.9998999999998

# Row 116
$ws.Range("H116").Value = 1290.375
$ws.Range("I116").Value = 1046.1428
$ws.Range("J116").Value = 3000
$ws.Range("K116").Value = 3138.4284
$ws.Range("L116").Value = 9000
$ws.Range("M116").Value = 303.5715999999998
$ws.Range("N116").Value = -15884

# Row 126
$ws.Range("H126").Value = 0
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 0
$ws.Range("L126").Value = 0
$ws.Range("M126").ClearContents()
$ws.Range("N126").ClearContents()

# Row 129
$ws.Range("H129").Value = 1260.6666
$ws.Range("J129").Value = 750
$ws.Range("L129").Value = 2250
$ws.Range("N129").Value = -12250

# Row 131
$ws.Range("H131").Value = 1580.6666
$ws.Range("I131").Value = 997.2
$ws.Range("K131").Value = 2991.6
$ws.Range("M131").Value = 2048.4

# Row 134
$ws.Range("H134").Value = 750
$ws.Range("I134").Value = 750
$ws.Range("K134").Value = 2250
$ws.Range("M134").Value = 2820

# Row 137
$ws.Range("H137").Value = 5400
$ws.Range("I137").Value = 800
$ws.Range("K137").Value = 2400
$ws.Range("M137").Value = 2700

# Row 138
$ws.Range("H138").Value = 6948.8335
$ws.Range("I138").Value = 5971.25
$ws.Range("K138").Value = 17913.75
$ws.Range("M138").Value = -12773.75

# Row 139
$ws.Range("H139").Value = 2208.7
$ws.Range("J139").Value = 1800
$ws.Range("L139").Value = 5400
$ws.Range("N139").Value = -15680

# Row 140
$ws.Range("H140").Value = 3414.6
$ws.Range("I140").Value = 2239.8
$ws.Range("J140").Value = 4589.4
$ws.Range("K140").Value = 6719.400000000001
$ws.Range("L140").Value = 13768.2
$ws.Range("M140").Value = -1539.400000000001
$ws.Range("N140").Value = -24128.2

$ws = $wb.Worksheets.Item("GSM")
# Row 3
$ws.Range("H3").Value = 1334.5454
$ws.Range("I3").Value = 279
$ws.Range("J3").Value = 3181.75
$ws.Range("K3").Value = 279
$ws.Range("L3").Value = 3181.75
$ws.Range("M3").Value = -163
$ws.Range("N3").Value = -3413.75

# Row 7
$ws.Range("H7").Value = 12501
$ws.Range("I7").Value = 10000.5
$ws.Range("J7").Value = 15001.5
$ws.Range("K7").Value = 10000.5
$ws.Range("L7").Value = 15001.5
$ws.Range("M7").Value = -9888.5
$ws.Range("N7").Value = -15225.5

# Row 8
$ws.Range("H8").Value = 12501
$ws.Range("I8").Value = 10000.5
$ws.Range("J8").Value = 15001.5
$ws.Range("K8").Value = 10000.5
$ws.Range("L8").Value = 15001.5
$ws.Range("M8").Value = -9861.5
$ws.Range("N8").Value = -15279.5

# Row 126
$ws.Range("H126").Value = 1955
$ws.Range("I126").Value = 1955
$ws.Range("K126").Value = 5865
$ws.Range("M126").Value = -3395

# Row 131
$ws.Range("H131").Value = 0
$ws.Range("J131").Value = 0
$ws.Range("L131").Value = 0
$ws.Range("N131").ClearContents()

# Row 132
$ws.Range("H132").Value = 3193.1428
$ws.Range("I132").Value = 2882.182
$ws.Range("J132").Value = 4333.3335
$ws.Range("K132").Value = 8646.545999999998
$ws.Range("L132").Value = 13000.0005
$ws.Range("M132").Value = -6116.545999999998
$ws.Range("N132").Value = -18060.0005

$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 3399.5
$ws.Range("I7").Value = 3399.5
$ws.Range("K7").Value = 3399.5
$ws.Range("M7").Value = -3287.5

# Row 22
$ws.Range("H22").Value = 1754.591
$ws.Range("I22").Value = 1818
$ws.Range("J22").Value = 1663
$ws.Range("K22").Value = 1818
$ws.Range("L22").Value = 1663
$ws.Range("M22").Value = -1523
$ws.Range("N22").Value = -2253

# Row 27
$ws.Range("H27").Value = 1754.591
$ws.Range("I27").Value = 1818
$ws.Range("J27").Value = 1663
$ws.Range("K27").Value = 1818
$ws.Range("L27").Value = 1663
$ws.Range("M27").Value = -1711
$ws.Range("N27").Value = -1877

# Row 126
$ws.Range("H126").Value = 3399.5
$ws.Range("I126").Value = 3399.5
$ws.Range("K126").Value = 10198.5
$ws.Range("M126").Value = -7728.5

$ws = $wb.Worksheets.Item("WVR")
# Row 6
$ws.Range("H6").Value = 252.5
$ws.Range("I6").Value = 252.5
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = 252.5
$ws.Range("L6").Value = 0
$ws.Range("M6").Value = -137.5
